# Realign ids names and titles for all new profiles for consistency
# (FHIR-36728) - update the US Core profile name column (B) for several
# rows, plus the FHIR Resource column (C) for the SDOH Assessment row,
# so the "US Core <Resource> <Topic> Profile" naming convention is used
# consistently, and reword the SDOH assessment profile/resource list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clinical Tests: Clinical Test / Clinical Test Result-Report
$ws.Cells.Item(23, 2).Value = "US Core Observation Clinical Test Result Profile, US Core DiagnosticReport Profile for Report and Note exchange"
$ws.Cells.Item(24, 2).Value = "US Core Observation Clinical Test Result Profile, US Core DiagnosticReport Profile for Report and Note exchange"

# Diagnostic Imaging: Diagnostic Imaging Test / Diagnostic Imaging Result-Report
$ws.Cells.Item(26, 2).Value = "US Core Observation Diagnostic Imaging Result Profile, US Core DiagnosticReport Profile for Report and Note exchange"
$ws.Cells.Item(27, 2).Value = "US Core Observation Diagnostic Imaging Result Profile, US Core DiagnosticReport Profile for Report and Note exchange"

# Encounter: Diagnosis
$ws.Cells.Item(30, 2).Value = "US Core Condition Encounter Diagnosis Profile"

# Goals: Health Concerns
$ws.Cells.Item(36, 2).Value = "US Core Condition Problems and Health Concerns Profile"

# Problems: header, Date of Resolution, Date of Diagnosis
$ws.Cells.Item(60, 2).Value = "US Core Condition Problems and Health Concerns Profile"
$ws.Cells.Item(61, 2).Value = "US Core Condition Problems and Health Concerns Profile"
$ws.Cells.Item(62, 2).Value = "US Core Condition Problems and Health Concerns Profile"

# Sexual Orientation
$ws.Cells.Item(67, 2).Value = "US Core Observation Sexual Orientation Profile"

# SDOH: Assessment
$ws.Cells.Item(69, 2).Value = "US Core Observation SDOH Assessment Profile,US Core Observation Social History Profile, US Core QuestionnaireResponse Profile"
$ws.Cells.Item(69, 3).Value = "Observation, QuestionnaireResponse"

# SDOH: Problems/Health Concerns
$ws.Cells.Item(72, 2).Value = "US Core Condition Problems and Health Concerns Profile"

# Cosmetic view-state touch-ups: zoom in to 130% and leave the selection
# parked on A9, matching the author's last-saved window state.
$excel.ActiveWindow.Zoom = 130
$ws.Range("A9").Select()
